$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1) Rename the sound-file text values used throughout the table ---
# Order matters: first move the existing "Coo_01.wav" entries out of the
# way (they become "Coo_02.wav"), THEN repoint the old "CooB.wav" entries
# onto "Coo_01.wav", and finally rename "CooA.wav" to "Avg16.wav".
$ws.Cells.Replace("Coo_01.wav", "Coo_02.wav")
$ws.Cells.Replace("CooB.wav", "Coo_01.wav")
$ws.Cells.Replace("CooA.wav", "Avg16.wav")

# --- 2) Extend the table with 4 more repeats of the 4-row cycle ---
# (rows 2-5 hold one full cycle incl. formatting; copy it down 4 times)
$cycle = $ws.Range("A2:C5")
$cycle.Copy($ws.Range("A50:C53"))
$cycle.Copy($ws.Range("A54:C57"))
$cycle.Copy($ws.Range("A58:C61"))
$cycle.Copy($ws.Range("A62:C65"))

# --- 3) Update the view/selection state ---
$excel.ActiveWindow.ScrollRow = 52
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A75").Select()
